# Update the organisation / contact-person information block on the
# "Пример" worksheet (rows 6-10, column B) with the newly supplied
# National Statistical Committee details.
#
# Row  6  (Организация)                        -> new department name
# Row  7  (Контактное лицо (лица) / Координатор) -> new contact person
# Row  8  (Электронная почта контактного лица)   -> new e-mail
# Row  9  (Телефон контактного лица)              -> new phone number
# Row 10  (Сайт организации (если есть))          -> new website
#
# The cells are unlocked on an otherwise protected sheet, so they can be
# edited directly without first calling Unprotect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = "Калымбетова Ы.И."
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B6").Value  = "Национальный статистический комитет (Управление статистики домашних хозяйств) в рамках глобальной программы MICS ЮНИСЕФ"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com"

# Leave the cursor on the last-touched cell, matching the saved selection.
$ws.Range("B10").Select() | Out-Null
